$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 91, shifting existing rows 91:141 down to 92:142
$ws.Rows("91:91").Insert()

# Populate the new row 91 with the new record's data.
# Columns A,B,C,E,F,G,H,N,O,Q,R are constant across this block of rows (same
# as row 90/92), so copy them straight from row 90.
$ws.Range("A91").Value = 10
$ws.Range("B91").Value = "Vega Modelo de Temuco"
$ws.Range("C91").Value = "La Araucanía"
$ws.Range("D91").Value = 44488
$ws.Range("D91").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E91").Value = 9
$ws.Range("F91").Value = 100112005
$ws.Range("G91").Value = "Puerro"
$ws.Range("H91").Value = "Azul de Maquehue"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 80
$ws.Range("K91").Value = 6000
$ws.Range("L91").Value = 7000
$ws.Range("M91").Value = 6500
$ws.Range("N91").Value = "`$/docena de paquetes"
$ws.Range("O91").Value = "Provincia de Cautín"
$ws.Range("P91").Value = 542
$ws.Range("Q91").Value = 12
$ws.Range("R91").Value = "Hortaliza"
